# Updates the cryptocurrency price / volume snapshot to the latest
# scrape (GitHub Actions scheduled refresh).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2-41: refresh Price (col D) / Volume(1h) (col E) text values ---
# Each entry is @(row, column, newText). NumberFormat is forced to "@"
# (Text) before the write so Excel does not reinterpret numeric-looking
# strings such as "22.491.88" or "0.3662" as numbers/dates.
$cellUpdates = @(
    @(2, 4, '22.491.88'),
    @(2, 5, '  +0.05%  '),
    @(3, 4, '1.573.33'),
    @(3, 5, '  +0.02%  '),
    @(4, 5, '  +0.05%  '),
    @(5, 5, '  +0.05%  '),
    @(6, 4, '287.16'),
    @(6, 5, '  -1.53%  '),
    @(7, 4, '0.3662'),
    @(7, 5, '  -1.58%  '),
    @(8, 4, '48.65'),
    @(8, 5, '  -2.56%  '),
    @(9, 4, '0.3340'),
    @(9, 5, '  -1.79%  '),
    @(10, 4, '1.132'),
    @(10, 5, '  -1.22%  '),
    @(11, 4, '0.07448'),
    @(11, 5, '  -1.40%  '),
    @(12, 5, '  +0.06%  '),
    @(13, 4, '20.86'),
    @(13, 5, '  -2.24%  '),
    @(14, 4, '5.994'),
    @(14, 5, '  -0.96%  '),
    @(15, 4, '6.933'),
    @(15, 5, '  -0.51%  '),
    @(16, 4, '1.577.24'),
    @(16, 5, '  +0.39%  '),
    @(17, 4, '0.00001112'),
    @(17, 5, '  -1.13%  '),
    @(18, 4, '88.34'),
    @(18, 5, '  -2.69%  '),
    @(19, 5, '  -0.04%  '),
    @(20, 5, '  +0.09%  '),
    @(21, 4, '6.399'),
    @(21, 5, '  +1.52%  '),
    @(22, 5, '  +0.33%  '),
    @(23, 5, '  -0.16%  '),
    @(24, 4, '22.481.66'),
    @(24, 5, '  -0.02%  '),
    @(25, 4, '2.386'),
    @(26, 4, '2.622'),
    @(26, 5, '  -0.74%  '),
    @(27, 4, '152.88'),
    @(27, 5, '  +2.32%  '),
    @(28, 4, '19.61'),
    @(28, 5, '  -2.17%  '),
    @(29, 4, '5.009'),
    @(29, 5, '  -0.88%  '),
    @(30, 4, '124.09'),
    @(30, 5, '  -0.91%  '),
    @(31, 4, '1.755.06'),
    @(31, 5, '  +0.43%  '),
    @(32, 4, '1.051'),
    @(32, 5, '  -3.26%  '),
    @(33, 4, '6.175'),
    @(33, 5, '  -1.03%  '),
    @(34, 4, '2.004'),
    @(34, 5, '  -0.24%  '),
    @(35, 4, '9.834'),
    @(35, 5, '  +0.46%  '),
    @(36, 4, '0.08298'),
    @(36, 5, '  -0.78%  '),
    @(37, 4, '0.02447'),
    @(37, 5, '  -1.63%  '),
    @(38, 4, '0.2266'),
    @(38, 5, '  -1.70%  '),
    @(39, 4, '0.06465'),
    @(39, 5, '  -0.88%  '),
    @(40, 4, '5.448'),
    @(40, 5, '  -0.44%  '),
    @(41, 4, '1.302'),
    @(41, 5, '  -2.59%  ')
)

foreach ($u in $cellUpdates) {
    $cell = $ws.Cells.Item($u[0], $u[1])
    $cell.NumberFormat = "@"
    $cell.Value = $u[2]
}

# --- Rows 42-51: coin ranking reshuffled ---
# TheSandbox/Aptos swapped places, "Frax" is newly inserted at row 44
# (pushing EnergySwap..Cronos down by one row), and "Aave" drops off the
# bottom of the list. Rewrite Coin (B), Link (C), Price (D) and
# Volume(1h) (E) for each of these rows in full.
$rowUpdates = @(
    @(42, 'TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.6362', '  +1.83%  '),
    @(43, 'Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '11.36', '  -0.14%  '),
    @(44, 'Frax', 'https://coinranking.com/coin/KfWtaeV1W+frax-frax', '1.002', '  +0.10%  '),
    @(45, 'EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '13.97', '  -0.88%  '),
    @(46, 'Decentraland', 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana', '0.6190', '  +5.23%  '),
    @(47, 'PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '3.764', '  -1.19%  '),
    @(48, 'NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '2.057', '  -1.08%  '),
    @(49, 'Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '125.40', '  -4.17%  '),
    @(50, 'EOS', 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos', '1.219', '  +0.27%  '),
    @(51, 'Cronos', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', '0.07250', '  -1.15%  ')
)

foreach ($u in $rowUpdates) {
    $row = $u[0]

    $coinCell = $ws.Cells.Item($row, 2)
    $coinCell.NumberFormat = "@"
    $coinCell.Value = $u[1]

    $linkCell = $ws.Cells.Item($row, 3)
    $linkCell.NumberFormat = "@"
    $linkCell.Value = $u[2]

    $priceCell = $ws.Cells.Item($row, 4)
    $priceCell.NumberFormat = "@"
    $priceCell.Value = $u[3]

    $volCell = $ws.Cells.Item($row, 5)
    $volCell.NumberFormat = "@"
    $volCell.Value = $u[4]
}

